$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy header style (bold, centered, bordered) from an existing header row
$weekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Copy date-format style from an existing date column
$weekly.Range("A2:A3").Copy()
$newSheet.Range("A2:A56").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Forecast data rows 2-56: ds, PO_Forecast, yhat_lower, yhat_upper ---
$data = @(
    @(44941.99999999999, 4, -3.459187432003717, 11.12956881244256),
    @(44948.99999999999, 4, -3.294752348062235, 11.43983420441185),
    @(44962.99999999999, 4, -2.945921955864102, 11.09129100981513),
    @(44969.99999999999, 4, -2.989628419641477, 11.21747408787788),
    @(44983.99999999999, 4, -2.726449097473561, 11.52007344242019),
    @(44990.99999999999, 4, -2.554476898335556, 11.62008645308434),
    @(44997.99999999999, 4, -3.058262558279143, 11.88049905904823),
    @(45004.99999999999, 4, -2.627401731060935, 11.86674081674074),
    @(45011.99999999999, 4, -3.008773813511833, 11.63870115376237),
    @(45018.99999999999, 4, -3.529538642707798, 11.3098307131423),
    @(45025.99999999999, 4, -3.415823176590115, 11.31631533160212),
    @(45039.99999999999, 4, -2.40932902622054, 11.55294942800163),
    @(45046.99999999999, 5, -2.597635489176695, 11.61845996846137),
    @(45053.99999999999, 5, -2.489042202053091, 11.82712275364829),
    @(45186.99999999999, 5, -1.614053688567775, 11.98073033090723),
    @(45193.99999999999, 5, -2.448208313567758, 12.95852102848749),
    @(45200.99999999999, 5, -1.830634819756759, 12.95705814512148),
    @(45207.99999999999, 5, -1.700145199293903, 12.86658455600246),
    @(45214.99999999999, 5, -2.27581041866378, 12.79988939072296),
    @(45221.99999999999, 6, -1.637802465392126, 13.13626050868089),
    @(45235.99999999999, 6, -1.936458005281587, 12.59690014620611),
    @(45270.99999999999, 6, -1.743418284336791, 12.97534733893017),
    @(45277.99999999999, 6, -1.539679330543523, 13.42203407466745),
    @(45333.99999999999, 6, -1.067450723686077, 13.35082628020225),
    @(45347.99999999999, 6, -1.066360201338182, 13.74002733412394),
    @(45368.99999999999, 6, -0.8916226556958714, 13.47915980972124),
    @(45375.99999999999, 6, -1.366276076412769, 13.50358543169474),
    @(45389.99999999999, 7, 0.07170871761246579, 13.60397514793496),
    @(45424.99999999999, 7, -0.4783039991759502, 13.88889677183683),
    @(45431.99999999999, 7, -0.1046113682831635, 14.3696623176038),
    @(45438.99999999999, 7, -0.4306927020226833, 13.79402517816193),
    @(45445.99999999999, 7, 0.2604303786226778, 14.19803149454809),
    @(45452.99999999999, 7, 0.1030664439329522, 14.05328549254824),
    @(45466.99999999999, 7, -0.2780615348225551, 13.74313589935023),
    @(45473.99999999999, 7, 0.07284798403708813, 14.83729530838492),
    @(45487.99999999999, 7, -0.5005727071411495, 14.09044427848382),
    @(45508.99999999999, 7, -0.08826194100228382, 14.548319409058),
    @(45522.99999999999, 7, -0.02147480892832063, 14.46786334137371),
    @(45529.99999999999, 7, 0.6550714514330404, 14.73939572444676),
    @(45536.99999999999, 7, 0.1731128109024191, 14.85084981532663),
    @(45557.99999999999, 7, 0.5470398973290324, 14.18432461532699),
    @(45564.99999999999, 8, 0.4520442082161557, 14.73783398206344),
    @(45606.99999999999, 8, 0.5809826383228893, 14.31380474522484),
    @(45613.99999999999, 8, 0.6381691834612138, 14.27144282467402),
    @(45627.99999999999, 8, 0.6418904792359234, 14.67596898254518),
    @(45634.99999999999, 8, 0.7102955160028205, 15.01489910741582),
    @(45641.99999999999, 8, 1.114562044179838, 15.07779599450731),
    @(45648.99999999999, 8, 0.711293896772732, 15.56763900211589),
    @(45655.99999999999, 8, 0.5447793647492124, 15.24560962905195),
    @(45662.99999999999, 8, 1.585377529939495, 15.54214420052116),
    @(45669.99999999999, 8, 0.4227153753205803, 15.26878266818178),
    @(45676.99999999999, 8, 0.9766992671039513, 15.07962174566011),
    @(45683.99999999999, 8, 1.300064819825062, 15.37726219480112),
    @(45690.99999999999, 8, 1.637929657974206, 16.12303324309471),
    @(45697.99999999999, 8, 1.101517754499828, 15.38550431577135)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
    $newSheet.Cells.Item($row, 3).Value = $data[$i][2]
    $newSheet.Cells.Item($row, 4).Value = $data[$i][3]
}

$newSheet.Range("A1").Select()